$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Aggregated Monthly Infographic")

# Update the D-column formulas (each is a string literal formula ="NNN")
$ws.Range("D2").Formula  = '="131"'
$ws.Range("D3").Formula  = '="127"'
$ws.Range("D4").Formula  = '="96.47"'
$ws.Range("D6").Formula  = '="0.05"'
$ws.Range("D7").Formula  = '="389"'
$ws.Range("D8").Formula  = '="2.17"'
$ws.Range("D9").Formula  = '="56"'
$ws.Range("D10").Formula = '="3.6"'
$ws.Range("D16").Formula = '="22"'
$ws.Range("D17").Formula = '="42"'

# Window position metadata change (xWindow/yWindow 1920 -> 1152)
$excel.ActiveWindow.Left = 1152
$excel.ActiveWindow.Top = 1152
